$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 136, shifting existing rows 136:163 down to 138:165
$ws.Rows("136:137").Insert()

# Row 136 (new)
$ws.Range("A136").Value = 1
$ws.Range("B136").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C136").Value = "Arica y Parinacota"
$ws.Range("D136").Value = 45005
$ws.Range("D136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E136").Value = 15
$ws.Range("F136").Value = 100112042
$ws.Range("G136").Value = "Locoto"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 40
$ws.Range("K136").Value = 45000
$ws.Range("L136").Value = 50000
$ws.Range("M136").Value = 47500
$ws.Range("N136").Value = "`$/caja 20 kilos"
$ws.Range("O136").Value = "Región de Arica y Parinacota"
$ws.Range("P136").Value = 2375
$ws.Range("Q136").Value = 20
$ws.Range("R136").Value = "Hortaliza"

# Row 137 (new)
$ws.Range("A137").Value = 1
$ws.Range("B137").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C137").Value = "Arica y Parinacota"
$ws.Range("D137").Value = 45005
$ws.Range("D137").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E137").Value = 15
$ws.Range("F137").Value = 100112042
$ws.Range("G137").Value = "Locoto"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 50
$ws.Range("K137").Value = 40000
$ws.Range("L137").Value = 45000
$ws.Range("M137").Value = 42500
$ws.Range("N137").Value = "`$/caja 20 kilos"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 2125
$ws.Range("Q137").Value = 20
$ws.Range("R137").Value = "Hortaliza"
